$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "2024-06-14 18:04:45"
$ws.Range("D4").Value = 200
$ws.Range("E4").Value = 1

# Row 5
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = "2024-06-14 18:04:46"
$ws.Range("D5").Value = 200
$ws.Range("E5").Value = 1
